$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (H1, bold+border+center/top alignment) onto new headers I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row data: row -> (I0 value, IF value)
$data = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(3, 4)
    5 = @(9, 9)
    6 = @(7, 8)
    7 = @(8, 8)
    8 = @(5, 5)
    9 = @(9, 9)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(5, 5)
    13 = @(7, 7)
    14 = @(6, 6)
    15 = @(9, 9)
    16 = @(7, 7)
    17 = @(6, 6)
    18 = @(6, 6)
    19 = @(6, 7)
    20 = @(6, 6)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(5, 5)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(6, 6)
    30 = @(8, 8)
    31 = @(7, 7)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(5, 6)
    35 = @(7, 7)
    36 = @(7, 7)
    37 = @(7, 8)
    38 = @(6, 7)
    39 = @(8, 9)
    40 = @(8, 8)
    41 = @(8, 9)
    42 = @(8, 8)
    43 = @(6, 6)
    44 = @(9, 9)
    45 = @(9, 9)
    46 = @(8, 8)
    47 = @(7, 7)
    48 = @(8, 8)
    49 = @(8, 8)
    50 = @(8, 8)
    51 = @(8, 8)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(8, 8)
    55 = @(7, 7)
    56 = @(7, 7)
    57 = @(7, 7)
    58 = @(8, 8)
    59 = @(8, 8)
    60 = @(8, 8)
    61 = @(8, 8)
    62 = @(7, 7)
    63 = @(9, 9)
    64 = @(8, 8)
    65 = @(8, 8)
    66 = @(8, 8)
    67 = @(8, 8)
    68 = @(9, 9)
    69 = @(8, 8)
    70 = @(7, 7)
    71 = @(8, 8)
    72 = @(8, 8)
    73 = @(9, 9)
    74 = @(8, 8)
    75 = @(6, 6)
    76 = @(7, 7)
    77 = @(7, 7)
    78 = @(4, 4)
    79 = @(3, 3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Host "done"